$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "exceluser2@test.com"
$ws.Range("B4").Value = "exceluser123123"
$ws.Range("A5").Value = "exceluser3@test.com"
$ws.Range("B5").Value = "exceluser123123"

$ws.Range("B5").Select()
